$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.722.51"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.649.82"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.536"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.19"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "1.882.44"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.636.99"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "27.695.36"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +9.67%  "
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").Value = "1.434.69"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.899"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.28%  "
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("B45").Value = "mCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "1.791.28"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("E50").Value = "  +7.41%  "
$ws.Range("E51").Value = "  -2.14%  "
